$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.494098544120789
$ws.Range("B1").Value = 2.052140951156616
$ws.Range("C1").Value = 2.417945861816406
$ws.Range("D1").Value = 2.887118816375732
$ws.Range("E1").Value = 2.618740558624268
